$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = 'Bitcoin'
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '98.562.42'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.32%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = 'Ethereum'
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").NumberFormat = "@"
$ws.Range("C3").Value = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range("C3").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.360.19'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.56%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = 'TetherUSD'
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = 'Solana'
$ws.Range("B5").Style = "Normal"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range("C5").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '256.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -0.85%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = 'BNB'
$ws.Range("B6").Style = "Normal"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '663.48'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.11%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("B7").NumberFormat = "@"
$ws.Range("B7").Value = 'XRP'
$ws.Range("B7").Style = "Normal"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range("C7").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.55'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +7.70%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("B8").NumberFormat = "@"
$ws.Range("B8").Value = 'Dogecoin'
$ws.Range("B8").Style = "Normal"
$ws.Range("C8").NumberFormat = "@"
$ws.Range("C8").Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range("C8").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.470'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +19.87%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = 'Cardano'
$ws.Range("B9").Style = "Normal"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range("C9").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '1.08'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +21.97%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = 'USDC'
$ws.Range("B10").Style = "Normal"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range("C10").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.999'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -0.03%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("B11").NumberFormat = "@"
$ws.Range("B11").Value = 'LidoStakedEther'
$ws.Range("B11").Style = "Normal"
$ws.Range("C11").NumberFormat = "@"
$ws.Range("C11").Value = 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
$ws.Range("C11").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '3.357.84'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +0.54%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("B12").NumberFormat = "@"
$ws.Range("B12").Value = 'TRON'
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").NumberFormat = "@"
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.216'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +8.57%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("B13").NumberFormat = "@"
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").NumberFormat = "@"
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '42.47'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +13.13%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000272'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +9.06%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'WrappedBTC'
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '98.320.46'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.24%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("B16").NumberFormat = "@"
$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.985.03'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +0.49%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("B17").NumberFormat = "@"
$ws.Range("B17").Value = 'Toncoin'
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").NumberFormat = "@"
$ws.Range("C17").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.68'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.81%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("B18").NumberFormat = "@"
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").NumberFormat = "@"
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.84'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +28.10%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.349.61'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.23%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'Chainlink'
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +10.64%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'BitcoinCash'
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '530.38'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.94%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("B22").NumberFormat = "@"
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("B22").Style = "Normal"
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("C22").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.56'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("B23").NumberFormat = "@"
$ws.Range("B23").Value = 'Uniswap'
$ws.Range("B23").Style = "Normal"
$ws.Range("C23").NumberFormat = "@"
$ws.Range("C23").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("C23").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.54'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +12.19%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("B24").NumberFormat = "@"
$ws.Range("B24").Value = 'PEPE'
$ws.Range("B24").Style = "Normal"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("C24").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000216'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.37%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("B25").NumberFormat = "@"
$ws.Range("B25").Value = 'Stellar'
$ws.Range("B25").Style = "Normal"
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("C25").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.443'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +50.11%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = 'Litecoin'
$ws.Range("B26").Style = "Normal"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("C26").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '102.80'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +14.97%  '
$ws.Range("E26").Style = "Normal"

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'NEARProtocol'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '6.25'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +10.90%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'Aptos'
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '12.56'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +5.24%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("B29").NumberFormat = "@"
$ws.Range("B29").Value = 'Hedera'
$ws.Range("B29").Style = "Normal"
$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.151'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +9.59%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").NumberFormat = "@"
$ws.Range("B30").Value = 'Dai'
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").NumberFormat = "@"
$ws.Range("C30").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.22%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").NumberFormat = "@"
$ws.Range("B31").Value = 'InternetComputer(DFINITY)'
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.17'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +14.63%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").NumberFormat = "@"
$ws.Range("B32").Value = 'Cronos'
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").NumberFormat = "@"
$ws.Range("C32").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.191'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -1.33%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'Binance-PegBSC-USD'
$ws.Range("B33").Style = "Normal"
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.00'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +0.29%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'EthereumClassic'
$ws.Range("B34").Style = "Normal"
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("C34").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '29.51'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +4.54%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("B35").NumberFormat = "@"
$ws.Range("B35").Value = 'PolygonEcosystemToken'
$ws.Range("B35").Style = "Normal"
$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.543'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +18.04%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("B36").NumberFormat = "@"
$ws.Range("B36").Value = 'RenderToken'
$ws.Range("B36").Style = "Normal"
$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '7.81'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +7.07%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("B37").NumberFormat = "@"
$ws.Range("B37").Value = 'PancakeSwap'
$ws.Range("B37").Style = "Normal"
$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.11'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +8.41%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.158'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +5.35%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = 'Bittensor'
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '524.11'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +4.70%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("B40").NumberFormat = "@"
$ws.Range("B40").Value = 'MantraDAO'
$ws.Range("B40").Style = "Normal"
$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.93'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +7.55%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = 'Fetch.AI'
$ws.Range("B41").Style = "Normal"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.34'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.10%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = 'WhiteBITCoin'
$ws.Range("B42").Style = "Normal"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '24.72'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.73%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("B43").NumberFormat = "@"
$ws.Range("B43").Value = 'VeChain'
$ws.Range("B43").Style = "Normal"
$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0434'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +33.19%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("B44").Style = "Normal"
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("C44").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.40'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +2.26%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("B45").Style = "Normal"
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("C45").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.827'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +5.90%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'USDe'
$ws.Range("B46").Style = "Normal"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("C46").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'Filecoin'
$ws.Range("B47").Style = "Normal"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("C47").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.18'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +11.56%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Stacks'
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.07'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +6.35%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'ImmutableX'
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.55'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +13.16%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'OKB'
$ws.Range("B50").Style = "Normal"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("C50").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '51.10'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +11.45%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Cosmos'
$ws.Range("B51").Style = "Normal"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("C51").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.83'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +14.64%  '
$ws.Range("E51").Style = "Normal"
